$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "[Lauren Chenarides](https://dataifa.github.io/difa-project/Leadership_team.html), [Drew Hanks](https://dataifa.github.io/difa-project/Leadership_team.html)"
$ws.Range("F4").Value = "[George Davis](https://dataifa.github.io/difa-project/Leadership_team.html), [Joe Cummins](https://www.josephrcummins.com/)"
$ws.Range("F6").Value = "[Drew Hanks](https://dataifa.github.io/difa-project/Leadership_team.html)"
$ws.Range("F7").Value = "[Lauren Chenarides](https://dataifa.github.io/difa-project/Leadership_team.html)"
$ws.Range("F8").Value = "[Amelia Finaret](https://dataifa.github.io/difa-project/Leadership_team.html)"
$ws.Range("F9").Value = "[Lauren Chenarides](https://dataifa.github.io/difa-project/Leadership_team.html), [Drew Hanks](https://dataifa.github.io/difa-project/Leadership_team.html)"

$ws.Range("G11").Select()
